# Speed Tracking.xlsx - stress testing re-run: refreshed Insertion/Fetch
# benchmark numbers after the Postgres autoincrement sequence was reset
# (batch-insert glitch) + nudge the window position / scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cosmetic window chrome -------------------------------------------------
# Book window moved on screen (workbookView xWindow/yWindow).
$win = $excel.ActiveWindow
$win.Left = 1540
$win.Top = 4660

# The frozen pane (first column frozen) was scrolled one column to the right
# (topLeftCell E1 -> F1); the split itself (xSplit=1) is unchanged.
$win.ScrollColumn = 6

# --- Re-run results ----------------------------------------------------------
# Row 2 (Mongo)
$ws.Range("G2").Value = 3.3071999999999999
$ws.Range("H2").Value = 2.3847999999999998
$ws.Range("K2").Value = 10.4026

# Row 3 (PostgreSQL listings)
$ws.Range("G3").Value = 3.649
$ws.Range("H3").Value = 3.9079999999999999
$ws.Range("K3").Value = 2.2360000000000002

# Row 4 (PostgreSQL bookedDates)
$ws.Range("G4").Value = 3.9660000000000002
$ws.Range("H4").Value = 4.3912000000000004
$ws.Range("K4").Value = 7.0288000000000004

# Row 5 (Postgres total)
$ws.Range("G5").Value = 7.6158000000000001
$ws.Range("H5").Value = 8.2992000000000008
$ws.Range("I5").Value = 4.3920000000000003
$ws.Range("J5").Value = 3.6814
$ws.Range("K5").Value = 9.2639999999999993

# Restore the original selection (K3) so the view isn't left showing the
# scroll-probe cell above.
$ws.Range("K3").Select()
